$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set B-column text cells FIRST, in the exact order needed so the
# shared-string table is interned/rebuilt in the same sequence as the target file. ---
$ws.Range("B2").Value = 'I agree with the news story. '
$ws.Range("B3").Value = 'I''m a robot connected to smart thermostats via the internet.'
$ws.Range("B5").Value = 'I have almost always experienced such issues in the homes where I was.'
$ws.Range("B8").Value = 'I also agree. '
$ws.Range("B9").Value = 'I can detect when Internet connectivity weakens. '
$ws.Range("B10").Value = 'I can alert people when they should check the Internet connectivity, before connection shuts down. '
$ws.Range("B11").Value = 'I have often used this alert function, there are many issues.'
$ws.Range("B12").Value = 'I see all of your points; I agree with the news story too.'
$ws.Range("B13").Value = 'I have temperature sensors to detect when a room is too hot or too cold. '
$ws.Range("B14").Value = 'I can fix it when a thermostat is not working correctly. I have always experienced temperature problems in the homes where I have been.'
$ws.Range("B15").Value = 'I have always experienced temperature problems in the homes where I have been.'
$ws.Range("B4").Value = 'I know a lot about the technology to evaluate its performance. '
$ws.Range("B6").Value = 'I diagree with the news story. '
$ws.Range("B7").Value = 'In the homes where I was, I have experienced none of the issues mentioned in the news story.'

# --- Now fill in the numeric A/C columns and header (order here does not matter) ---
$ws.Range("C1").Value = "time"
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 0.3
$ws.Range("A3").Value = 1
$ws.Range("C3").Value = 0.3
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 0.3
$ws.Range("A5").Value = 1
$ws.Range("C5").Value = 1.5
$ws.Range("A6").Value = 4
$ws.Range("C6").Value = 0.3
$ws.Range("A7").Value = 4
$ws.Range("C7").Value = 1.5
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = 0.3
$ws.Range("A9").Value = 2
$ws.Range("C9").Value = 0.3
$ws.Range("A10").Value = 2
$ws.Range("C10").Value = 0.3
$ws.Range("A11").Value = 2
$ws.Range("C11").Value = 1.5
$ws.Range("A12").Value = 3
$ws.Range("C12").Value = 0.3
$ws.Range("A13").Value = 3
$ws.Range("C13").Value = 0.3
$ws.Range("A14").Value = 3
$ws.Range("C14").Value = 0.3
$ws.Range("A15").Value = 3

# --- View changes: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("B24").Select()
